$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date column (A) uses a bold/centered/bordered date-number style.
# Copy that cell's formatting down onto the two new date cells before
# filling in the values, so A64/A65 pick up the same style as A63.
$ws.Range("A63").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("A65").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 64 - 2024-07-21
$ws.Range("A64").Value = 45494
$ws.Range("B64").Value = 727.5450335016
$ws.Range("C64").Value = 245.167786692
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 130.1640504
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 270.5838690977
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 20.9967939854805
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("N64").Value = 163.31844516
$ws.Range("O64").Value = 61.199160425
$ws.Range("P64").Value = 0
$ws.Range("Q64").Value = 0.000003
$ws.Range("R64").Value = 0
$ws.Range("S64").Value = 0
$ws.Range("T64").Value = 0
$ws.Range("U64").Value = 388.3924990490556
$ws.Range("V64").Value = 0
$ws.Range("W64").Value = 0
$ws.Range("X64").Value = 0
$ws.Range("Y64").Value = 0
$ws.Range("Z64").Value = 288.68784614925

# Row 65 - 2024-07-22 (columns H and V are left blank, matching the source)
$ws.Range("A65").Value = 45495
$ws.Range("B65").Value = 720.7853504124
$ws.Range("C65").Value = 238.48930946
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 120.86019735
$ws.Range("I65").Value = 262.392370595
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 20.699290646328
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("N65").Value = 144.18924881408
$ws.Range("O65").Value = 59.645057264
$ws.Range("P65").Value = 0
$ws.Range("Q65").Value = 0.000002856
$ws.Range("R65").Value = 0
$ws.Range("S65").Value = 0
$ws.Range("T65").Value = 0
$ws.Range("U65").Value = 367.4121400753912
$ws.Range("W65").Value = 0
$ws.Range("X65").Value = 0
$ws.Range("Y65").Value = 0
$ws.Range("Z65").Value = 268.047965545282
